$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Sandeep"
$ws.Range("B2").Value = "Sandeep2@gmail.com"
$ws.Range("C2").Value = 7817008526
$ws.Range("D2").Value = "Vadodara"
$ws.Range("F2").Value = "sand"
$ws.Range("G2").Value = "sand"
